$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank cells on row 3
$ws.Cells.Item(3, 2).Value = "programacion"
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = "impacto al cliente"
$ws.Cells.Item(3, 5).Value = "grande empresa"

# Add new row 6: collection "e"
$ws.Cells.Item(6, 1).Value = "e"
$ws.Cells.Item(6, 2).Value = "programacion"
$ws.Cells.Item(6, 3).Value = 12
$ws.Cells.Item(6, 4).Value = "impacto al cliente"
$ws.Cells.Item(6, 5).Value = "pequeña empresa"
$ws.Cells.Item(6, 6).Value = "Usuario1"

# Add new row 7: collection "f"
$ws.Cells.Item(7, 1).Value = "f"
$ws.Cells.Item(7, 2).Value = "programacion"
$ws.Cells.Item(7, 3).Value = 15
$ws.Cells.Item(7, 4).Value = "impacto al cliente"
$ws.Cells.Item(7, 5).Value = "grande empresa"
$ws.Cells.Item(7, 6).Value = "Usuario2"

# Match the formatting of existing data rows (style index 2 / xf applied to rows 2-5)
$srcRow = $ws.Range("A5:F5")
$newRows = $ws.Range("A6:F7")
$srcRow.Copy()
$newRows.PasteSpecial(-4122)
$excel.CutCopyMode = $false
